$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. The US data is no longer used anywhere - remove the whole tab
#    first so every sheet reference grabbed afterwards reflects the
#    final tab order/indices.
# ---------------------------------------------------------------------
$wb.Worksheets.Item("Data_US").Delete() | Out-Null

# ---------------------------------------------------------------------
# 2. Update Data_Brazil: drop the old "fraction" scratch formula in C11,
#    add three new labelled fraction rows (14-16) that replace what the
#    FoFObE tab used to compute directly off Data_US.
# ---------------------------------------------------------------------
$brazil = $wb.Worksheets.Item("Data_Brazil")

$brazil.Range("C11").ClearContents() | Out-Null

$brazil.Range("A14").Value = "Domestic Industry Fraction"
$brazil.Range("A16").Value = "Government Fraction"
$brazil.Range("A15").Value = "Labour Fraction"

$brazil.Range("B14").Formula = "=(E9+K9)/B12"
$brazil.Range("B14").ClearFormats() | Out-Null

$brazil.Range("B15").Formula = "=(H9+N9)/B12"
$brazil.Range("B15").ClearFormats() | Out-Null

$brazil.Range("B16").Formula = "=B9/B12"
$brazil.Range("B16").ClearFormats() | Out-Null

$brazil.Range("C16").Select() | Out-Null

# ---------------------------------------------------------------------
# 3. Update FoFObE so it reads straight from Data_Brazil instead of the
#    US-report computations.
# ---------------------------------------------------------------------
$fofobe = $wb.Worksheets.Item("FoFObE")

$fofobe.Range("B2").Formula = "=Data_Brazil!B16"
$fofobe.Range("B3").Formula = "=Data_Brazil!B14"
$fofobe.Range("B4").Formula = "=Data_Brazil!B15"

$fofobe.Range("B5").Select() | Out-Null

# ---------------------------------------------------------------------
# 4. FoFObE (the "blue tab") is now the one analysts should land on.
# ---------------------------------------------------------------------
$fofobe.Activate()

$wb.Save()
